$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - new data row (product line)
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "CETAL 120 MG 5 INFANTILE SUPP"
$ws.Range("H4").Value = "8:0"
$ws.Range("L4").Value = 13
$ws.Range("N4").Value = "1:0"

# The "name" (B4:G4 merged) and "transactions" (N4) cells share one style;
# the "balance" (H4:K4 merged) cells share another. Both switch from the
# General number format to Text ("@" = numFmtId 49) now that they hold text.
$ws.Range("B4:G4").NumberFormat = "@"
$ws.Range("N4").NumberFormat = "@"
$ws.Range("H4:K4").NumberFormat = "@"

# Row 5 - totals row picks up the price total and grows slightly taller
$ws.Rows.Item(5).RowHeight = 26.25
$ws.Range("K5").Value = 13
